$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns (C:F) before the existing ExpPoints column, shifting
# ExpPoints from C to G. The inserted columns inherit the header row's style.
$ws.Range("C1:F1").EntireColumn.Insert()

# New header row
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# Placeholder (empty text) cells for the new WIN/TOP4/TOP5/RELEGATION columns,
# rows 2-21. A leading apostrophe forces a literal (non-formula) empty text
# cell instead of clearing the cell outright, then ClearFormats drops the
# implicit "quote prefix" style so the cell stays on the default style.
$placeholders = $ws.Range("C2:F21")
$placeholders.Value = "'"
$placeholders.ClearFormats()

# Updated team order + recomputed ExpPoints values (rank column A is unchanged)
$ws.Range("B2").Value = "Arsenal"
$ws.Range("G2").Value = 82.35203596254101

$ws.Range("B3").Value = "Manchester City"
$ws.Range("G3").Value = 72.5929176890066

$ws.Range("B4").Value = "Liverpool"
$ws.Range("G4").Value = 69.57884018614058

$ws.Range("B5").Value = "Chelsea"
$ws.Range("G5").Value = 62.29344302794853

$ws.Range("B6").Value = "Crystal Palace"
$ws.Range("G6").Value = 59.17888377607198

$ws.Range("B7").Value = "Aston Villa"
$ws.Range("G7").Value = 59.07763905117865

$ws.Range("B8").Value = "AFC Bournemouth"
$ws.Range("G8").Value = 56.16657373935881

$ws.Range("B9").Value = "Brighton & Hove Albion"
$ws.Range("G9").Value = 55.30100272275261

$ws.Range("B10").Value = "Manchester United"
$ws.Range("G10").Value = 54.13324275781807

$ws.Range("B11").Value = "Newcastle United"
$ws.Range("G11").Value = 53.84894638489681

$ws.Range("B12").Value = "Tottenham Hotspur"
$ws.Range("G12").Value = 52.7213371054105

$ws.Range("B13").Value = "Brentford"
$ws.Range("G13").Value = 51.59970338152792

$ws.Range("B14").Value = "Sunderland"
$ws.Range("G14").Value = 43.80177441763223

$ws.Range("B15").Value = "Everton"
$ws.Range("G15").Value = 43.74888897103195

$ws.Range("B16").Value = "Fulham"
$ws.Range("G16").Value = 43.23384846378813

$ws.Range("B17").Value = "Leeds United"
$ws.Range("G17").Value = 36.79862195775968

$ws.Range("B18").Value = "Nottingham Forest"
$ws.Range("G18").Value = 36.19368624372028

$ws.Range("B19").Value = "Burnley"
$ws.Range("G19").Value = 35.39473326640184

$ws.Range("B20").Value = "West Ham United"
$ws.Range("G20").Value = 35.32089427355386

$ws.Range("B21").Value = "Wolverhampton Wanderers"
$ws.Range("G21").Value = 27.54304857240759
